# Auto-generated script to apply market-price/profit recalculation updates
# to the Asura_Profits workbook (FFXIV leve profit tracker).
# For each affected row, columns H-N (computed market price / profit columns)
# are refreshed to the latest values pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 4124.625
$ws.Range("I32").Value = 6698.3335
$ws.Range("J32").Value = 2580.4
$ws.Range("K32").Value = 6698.3335
$ws.Range("L32").Value = 2580.4
$ws.Range("M32").Value = -6372.3335
$ws.Range("N32").Value = -3232.4
# Row 80
$ws.Range("H80").Value = 6851.1113
$ws.Range("I80").Value = 901.7778
$ws.Range("J80").Value = 12800.444
$ws.Range("K80").Value = 2705.3334
$ws.Range("L80").Value = 38401.33199999999
$ws.Range("M80").Value = -1707.3334
$ws.Range("N80").Value = -40397.33199999999
# Row 83
$ws.Range("H83").Value = 6851.1113
$ws.Range("I83").Value = 901.7778
$ws.Range("J83").Value = 12800.444
$ws.Range("K83").Value = 8116.000199999999
$ws.Range("L83").Value = 115203.996
$ws.Range("M83").Value = -3124.000199999999
$ws.Range("N83").Value = -125187.996
# Row 92
$ws.Range("H92").Value = 1011.05
$ws.Range("I92").Value = 1011.05
$ws.Range("K92").Value = 1011.05
$ws.Range("M92").Value = 236.95
# Row 112
$ws.Range("H112").Value = 2718.3
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").Value = $null
# Row 121
$ws.Range("H121").Value = 1707.2727
$ws.Range("J121").Value = 2128.5715
$ws.Range("L121").Value = 6385.7145
$ws.Range("N121").Value = -9879.7145
# Row 138
$ws.Range("H138").Value = 2350.2222
$ws.Range("I138").Value = 1314.5814
$ws.Range("J138").Value = 4576.85
$ws.Range("K138").Value = 3943.7442
$ws.Range("L138").Value = 13730.55
$ws.Range("M138").Value = 1196.2558
$ws.Range("N138").Value = -24010.55

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15989.667
$ws.Range("I32").Value = 15981.515
$ws.Range("K32").Value = 15981.515
$ws.Range("M32").Value = -15694.515
# Row 45
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -2254
# Row 61
$ws.Range("H61").Value = 3220.9167
$ws.Range("I61").Value = 3228.5454
$ws.Range("J61").Value = 3214.4614
$ws.Range("K61").Value = 3228.5454
$ws.Range("L61").Value = 3214.4614
$ws.Range("M61").Value = -3016.5454
$ws.Range("N61").Value = -3638.4614
# Row 92
$ws.Range("H92").Value = 72595.55499999999
$ws.Range("J92").Value = 72595.55499999999
$ws.Range("L92").Value = 72595.55499999999
$ws.Range("N92").Value = -77587.55499999999
# Row 97
$ws.Range("H97").Value = 987.1111
$ws.Range("I97").Value = 681.0714
$ws.Range("K97").Value = 681.0714
$ws.Range("M97").Value = -185.0714
# Row 122
$ws.Range("H122").Value = 6422.6
$ws.Range("I122").Value = 8287.429
$ws.Range("J122").Value = 4790.875
$ws.Range("K122").Value = 24862.287
$ws.Range("L122").Value = 14372.625
$ws.Range("M122").Value = -22412.287
$ws.Range("N122").Value = -19272.625
# Row 136
$ws.Range("H136").Value = 3220.9167
$ws.Range("I136").Value = 3228.5454
$ws.Range("J136").Value = 3214.4614
$ws.Range("K136").Value = 9685.636200000001
$ws.Range("L136").Value = 9643.3842
$ws.Range("M136").Value = -7135.636200000001
$ws.Range("N136").Value = -14743.3842

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 52609.9
$ws.Range("I20").Value = 84844.836
$ws.Range("J20").Value = 4257.5
$ws.Range("K20").Value = 84844.836
$ws.Range("L20").Value = 4257.5
$ws.Range("M20").Value = -84597.836
$ws.Range("N20").Value = -4751.5
# Row 64
$ws.Range("H64").Value = 1347.8334
$ws.Range("I64").Value = 300
$ws.Range("J64").Value = 1871.75
$ws.Range("K64").Value = 300
$ws.Range("L64").Value = 1871.75
$ws.Range("M64").Value = -75
$ws.Range("N64").Value = -2321.75
# Row 67
$ws.Range("H67").Value = 1347.8334
$ws.Range("I67").Value = 300
$ws.Range("J67").Value = 1871.75
$ws.Range("K67").Value = 300
$ws.Range("L67").Value = 1871.75
$ws.Range("M67").Value = 480
$ws.Range("N67").Value = -3431.75
# Row 86
$ws.Range("H86").Value = 252862.75
$ws.Range("I86").Value = 4348.75
$ws.Range("J86").Value = 501376.75
$ws.Range("K86").Value = 4348.75
$ws.Range("L86").Value = 501376.75
$ws.Range("M86").Value = -3225.75
$ws.Range("N86").Value = -503622.75
# Row 89
$ws.Range("H89").Value = 252862.75
$ws.Range("I89").Value = 4348.75
$ws.Range("J89").Value = 501376.75
$ws.Range("K89").Value = 21743.75
$ws.Range("L89").Value = 2506883.75
$ws.Range("M89").Value = -16127.75
$ws.Range("N89").Value = -2518115.75
# Row 105
$ws.Range("H105").Value = 2685.077
$ws.Range("I105").Value = 2556
$ws.Range("J105").Value = 3395
$ws.Range("K105").Value = 2556
$ws.Range("L105").Value = 3395
$ws.Range("M105").Value = -809
$ws.Range("N105").Value = -6889
# Row 134
$ws.Range("H134").Value = 2523.7026
$ws.Range("I134").Value = 2208.9355
$ws.Range("K134").Value = 6626.806500000001
$ws.Range("M134").Value = -4091.806500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 248.40909
$ws.Range("I23").Value = 46.6
$ws.Range("J23").Value = 307.7647
$ws.Range("K23").Value = 139.8
$ws.Range("L23").Value = 923.2941000000001
$ws.Range("M23").Value = 95.19999999999999
$ws.Range("N23").Value = -1393.2941
# Row 62
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10372
# Row 65
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 27000
$ws.Range("N65").Value = -33864
# Row 113
$ws.Range("H113").Value = 545.125
$ws.Range("I113").Value = 556.05884
$ws.Range("J113").Value = 518.5714
$ws.Range("K113").Value = 1668.17652
$ws.Range("L113").Value = 1555.7142
$ws.Range("M113").Value = 501.82348
$ws.Range("N113").Value = -5895.7142
# Row 136
$ws.Range("H136").Value = 5427.533
$ws.Range("I136").Value = 748
$ws.Range("J136").Value = 7767.3
$ws.Range("K136").Value = 2244
$ws.Range("L136").Value = 23301.9
$ws.Range("M136").Value = 2856
$ws.Range("N136").Value = -33501.9

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8091.727
$ws.Range("I70").Value = 8000
$ws.Range("J70").Value = 8144.143
$ws.Range("K70").Value = 8000
$ws.Range("L70").Value = 8144.143
$ws.Range("M70").Value = -7730
$ws.Range("N70").Value = -8684.143
# Row 73
$ws.Range("H73").Value = 8091.727
$ws.Range("I73").Value = 8000
$ws.Range("J73").Value = 8144.143
$ws.Range("K73").Value = 8000
$ws.Range("L73").Value = 8144.143
$ws.Range("M73").Value = -7064
$ws.Range("N73").Value = -10016.143
# Row 102
$ws.Range("H102").Value = 2733.5
$ws.Range("I102").Value = 2302.6924
$ws.Range("K102").Value = 2302.6924
$ws.Range("M102").Value = -680.6923999999999
# Row 126
$ws.Range("H126").Value = 2559.5
$ws.Range("I126").Value = 1650
$ws.Range("J126").Value = 3014.25
$ws.Range("K126").Value = 4950
$ws.Range("L126").Value = 9042.75
$ws.Range("M126").Value = -2480
$ws.Range("N126").Value = -13982.75
# Row 136
$ws.Range("H136").Value = 42883.668
$ws.Range("J136").Value = 42883.668
$ws.Range("L136").Value = 128651.004
$ws.Range("N136").Value = -133751.004

$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
# Row 96
$ws.Range("H96").Value = 1876
$ws.Range("I96").Value = 1314.5
$ws.Range("K96").Value = 1314.5
$ws.Range("M96").Value = 58.5
